# inventory, loadfile goods, delete table invent group
#
# Adds three new inventory items below the existing "книга" row:
#   - row 3 (E3): "пираты"  -> "пираты2"   (re-label existing "Группа" value)
#   - row 4      : "тетрадь"  / Группа "пираты3"
#   - row 5      : "карандаш" / Группа "пираты3"
# and moves the active selection down to E6 (next empty "Группа" cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the group on the existing third data row.
$ws.Range("E3").Value = "пираты2"

# New row: тетрадь / пираты3
# Seed formatting for the "Группа" cell (E4) from the cell above it (E3),
# which carries the correct style, then overwrite the values.
$ws.Range("E3").Copy($ws.Range("E4"))
$ws.Range("A4").Value = "тетрадь"
$ws.Range("E4").Value = "пираты3"

# New row: карандаш / пираты3
$ws.Range("E3").Copy($ws.Range("E5"))
$ws.Range("A5").Value = "карандаш"
$ws.Range("E5").Value = "пираты3"

# Move the active cell to the next free "Группа" cell, as in the source edit.
$ws.Range("E6").Select()
